$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A13").Value = "24/10/2025"
$ws.Range("B13").Value = "Preston"
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = "Sheff Utd"
$ws.Range("F13").Value = "W"
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 2
$ws.Range("I13").Value = 2
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 1.57
$ws.Range("L13").Value = 2.37
$ws.Range("M13").Value = 15
$ws.Range("N13").Value = 12
$ws.Range("O13").Value = 3
$ws.Range("P13").Value = 3
